$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header: Volume/Number and date range (shared-string runs) ---
$volChars = $ws.Range("A8").Characters(21, 2)
$volChars.Text = "22"
$dateChars1 = $ws.Range("C9").Characters(27, 9)
$dateChars1.Text = "5/26/2025"
$dateChars2 = $ws.Range("C9").Characters(47, 9)
$dateChars2.Text = "6/1/2025"

# --- Update crime statistics table (rows 15-31) ---
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -60
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = -54.166666666666
$ws.Range("L16").Value = -56
$ws.Range("M16").Value = -69.444444444444
$ws.Range("N16").Value = -92.086330935251
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = 34.210526315789
$ws.Range("L17").Value = 13.333333333333
$ws.Range("M17").Value = 131.818181818182
$ws.Range("N17").Value = -5.555555555555
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 13.636363636363
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = 12.844036697247
$ws.Range("L18").Value = -15.753424657534
$ws.Range("M18").Value = 17.142857142857
$ws.Range("N18").Value = -71.981776765375
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 15.384615384615
$ws.Range("I19").Value = 174
$ws.Range("J19").Value = 197
$ws.Range("K19").Value = -11.675126903553
$ws.Range("L19").Value = -40
$ws.Range("M19").Value = 14.473684210526
$ws.Range("N19").Value = -17.535545023696
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -27.272727272727
$ws.Range("I20").Value = 101
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = -7.339449541284
$ws.Range("L20").Value = 62.903225806451
$ws.Range("M20").Value = 65.573770491803
$ws.Range("N20").Value = -92.417417417417
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 58.823529411764
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -3.529411764705
$ws.Range("I21").Value = 465
$ws.Range("J21").Value = 480
$ws.Range("K21").Value = -3.125
$ws.Range("L21").Value = -18.989547038327
$ws.Range("M21").Value = 22.691292875989
$ws.Range("N21").Value = -78.689275893675
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -41.176470588235
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 219
$ws.Range("J24").Value = 235
$ws.Range("K24").Value = -6.808510638297
$ws.Range("L24").Value = -1.793721973094
$ws.Range("M24").Value = 9.5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 5
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 43
$ws.Range("J25").Value = 48
$ws.Range("K25").Value = -10.416666666666
$ws.Range("L25").Value = -10.416666666666
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 300
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = -13.636363636363
$ws.Range("I26").Value = 81
$ws.Range("J26").Value = 83
$ws.Range("K26").Value = -2.409638554216
$ws.Range("L26").Value = -12.903225806451
$ws.Range("M26").Value = 19.117647058823
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Value = 0
$ws.Range("L14").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 0
$ws.Range("F28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -75
$ws.Range("N29").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("F31").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = 66.666666666666

$excel.CutCopyMode = 0